$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 15,
# shifting all subsequent rows (15-125) down by one (to 16-126).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with its own data.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44649
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112001
$ws.Range("G15").Value = "Berenjena"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8500
$ws.Range("N15").Value = "$/caja 50 unidades"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 170
$ws.Range("Q15").Value = 50
$ws.Range("R15").Value = "Hortaliza"
